# Pushing June 11 files
# Append four new pitchers to the headshot lookup table (rows 202-205).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Name) -------------------------------------------------
$ws.Range("A202").Value = "Drew Thorpe"
$ws.Range("A203").Value = "Jose Suarez"
$ws.Range("A204").Value = "Louie Varland"
$ws.Range("A205").Value = "Carlos Rodriguez"

# --- Column B (Baseball_Savant_Name) ----------------------------------
$ws.Range("B202").Value = "Drew Thorpe"
$ws.Range("B203").Value = "José Suarez"
$ws.Range("B204").Value = "Louie Varland"
$ws.Range("B205").Value = "Carlos Rodriguez"

# --- Column C (headshot URL / hyperlink) ------------------------------
$ws.Hyperlinks.Add($ws.Cells.Item(202, 3), "https://midfield.mlbstatic.com/v1/people/657514/spots/120")
$ws.Range("C202").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Cells.Item(203, 3), "https://midfield.mlbstatic.com/v1/people/660761/spots/120")
$ws.Range("C203").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Cells.Item(204, 3), "https://midfield.mlbstatic.com/v1/people/686973/spots/120")
$ws.Range("C204").Style = "Hyperlink"

# row 205 has no headshot URL / hyperlink

# --- Column D (Handedness) -------------------------------------------
$ws.Range("D202").Value = "RHP"
$ws.Range("D203").Value = "LHP"
$ws.Range("D204").Value = "RHP"
$ws.Range("D205").Value = "RHP"

# Leave the selection where the user would land after typing the new rows.
[void]$ws.Range("D206").Select()
